# Scheduled runner update: refresh market-board derived profit figures
# across the Pandaemonium_Profits leve sheets (ALC, ARM, BSM, CRP, CUL,
# GSM, LTW, WVR). Only the numeric currentAveragePrice* / LevePrice* /
# LeveProfit* columns (H:N) are touched; columns A:G are left untouched.

$wb = $excel.ActiveWorkbook

function Set-LeveRow {
    param(
        [string]$SheetName,
        [int]$Row,
        [hashtable]$Values,
        [bool]$RemoveN
    )

    $ws = $wb.Worksheets.Item($SheetName)

    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }

    if ($RemoveN) {
        $ws.Range("N$Row").ClearContents()
    }
}

# ---------------- ALC ----------------
Set-LeveRow "ALC" 86 @{
    H = 4355.8; I = 3444.75; J = 8000; K = 3444.75; L = 8000; M = -2321.75; N = -10246
} $false
Set-LeveRow "ALC" 89 @{
    H = 4355.8; I = 3444.75; J = 8000; K = 17223.75; L = 40000; M = -11607.75; N = -51232
} $false
Set-LeveRow "ALC" 112 @{
    H = 5478.8823; J = 1774.3334; L = 5323.0002; N = -7539.0002
} $false

# ---------------- ARM ----------------
Set-LeveRow "ARM" 74 @{
    H = 2183.8438; I = 1983.4231; J = 3052.3333; K = 1983.4231; L = 3052.3333; M = -1109.4231; N = -4800.3333
} $false
Set-LeveRow "ARM" 77 @{
    H = 2183.8438; I = 1983.4231; J = 3052.3333; K = 9917.1155; L = 15261.6665; M = -5549.1155; N = -23997.6665
} $false
Set-LeveRow "ARM" 88 @{
    H = 6199.5; J = 3038.7; L = 3038.7; N = -3850.7
} $false
Set-LeveRow "ARM" 91 @{
    H = 6199.5; J = 3038.7; L = 3038.7; N = -5846.7
} $false
Set-LeveRow "ARM" 132 @{
    H = 5140.931; I = 5753.409; K = 17260.227; M = -14730.227
} $false

# ---------------- BSM ----------------
Set-LeveRow "BSM" 20 @{
    H = 4166; I = 4603.2; J = 1980; K = 4603.2; L = 1980; M = -4356.2; N = -2474
} $false
Set-LeveRow "BSM" 86 @{
    H = 1833.9056; I = 1820.6123; J = 1996.75; K = 1820.6123; L = 1996.75; M = -697.6123; N = -4242.75
} $false
Set-LeveRow "BSM" 89 @{
    H = 1833.9056; I = 1820.6123; J = 1996.75; K = 9103.0615; L = 9983.75; M = -3487.0615; N = -21215.75
} $false
Set-LeveRow "BSM" 134 @{
    H = 2793.0645; I = 2441.0417; K = 7323.1251; M = -4788.1251
} $false

# ---------------- CRP ----------------
# Rows 20/30/128 lose their N cell entirely (no LeveProfitHQ value remains).
Set-LeveRow "CRP" 20 @{
    H = 0; I = 0; J = 0; K = 0; L = 0
} $true
Set-LeveRow "CRP" 30 @{
    H = 0; I = 0; J = 0; K = 0; L = 0
} $true
Set-LeveRow "CRP" 31 @{
    H = 6874.4653; I = 6626.528; J = 8149.5713; K = 6626.528; L = 8149.5713; M = -6331.528; N = -8739.5713
} $false
Set-LeveRow "CRP" 34 @{
    H = 6874.4653; I = 6626.528; J = 8149.5713; K = 6626.528; L = 8149.5713; M = -6424.528; N = -8553.5713
} $false
Set-LeveRow "CRP" 94 @{
    H = 1135.3077; I = 1212; J = 1128.9166; K = 1212; L = 1128.9166; M = -761; N = -2030.9166
} $false
Set-LeveRow "CRP" 99 @{
    H = 2220; I = 1550; J = 4900; K = 1550; L = 4900; M = -52; N = -7896
} $false
Set-LeveRow "CRP" 126 @{
    H = 2220; I = 1550; J = 4900; K = 4650; L = 14700; M = -2180; N = -19640
} $false
Set-LeveRow "CRP" 128 @{
    H = 0; I = 0; J = 0; K = 0; L = 0
} $true

# ---------------- CUL ----------------
Set-LeveRow "CUL" 34 @{
    H = 1999.174; J = 2720.9375; L = 8162.8125; N = -8330.8125
} $false
Set-LeveRow "CUL" 39 @{
    H = 5401; J = 5401; L = 16203; N = -16791
} $false
Set-LeveRow "CUL" 122 @{
    H = 942.0909; I = 232; J = 1099.8889; K = 2088; L = 9899.0001; M = 362; N = -14799.0001
} $false

# ---------------- GSM ----------------
Set-LeveRow "GSM" 70 @{
    H = 5642.4; I = 5476; J = 6100; K = 5476; L = 6100; M = -5206; N = -6640
} $false
Set-LeveRow "GSM" 73 @{
    H = 5642.4; I = 5476; J = 6100; K = 5476; L = 6100; M = -4540; N = -7972
} $false

# ---------------- LTW ----------------
Set-LeveRow "LTW" 136 @{
    H = 6449.9653; I = 3842.2307; J = 8568.75; K = 11526.6921; L = 25706.25; M = -8976.6921; N = -30806.25
} $false

# ---------------- WVR ----------------
Set-LeveRow "WVR" 113 @{
    H = 706.24; I = 414.15384; K = 1242.46152; M = 927.53848
} $false
Set-LeveRow "WVR" 132 @{
    H = 3279.7144; I = 2993.1667; K = 8979.5001; M = -6449.5001
} $false
